# Daily update at 8 AM UTC
# Appends the next day's row (row 66) to the Wins Over Time tracking sheet
# and re-applies the "last row" date number format to the new final row,
# restoring the previous last row (65) back to the standard date format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the special "last row" date format currently on A65 before we
# overwrite it, so it can be moved onto the new last row (A66).
$lastRowNumberFormat = $ws.Range("A65").NumberFormat

# Row 65 is no longer the last row -- give it the regular date format (same
# format used by all the other data rows, e.g. row 64).
$ws.Range("A65").NumberFormat = $ws.Range("A64").NumberFormat

# Add the new day's data in row 66, applying the "last row" date format that
# row 65 used to have onto the new final row.
$ws.Range("A66").NumberFormat = $lastRowNumberFormat
$ws.Range("A66").Value = 45806
$ws.Range("B66").Value = 279
$ws.Range("C66").Value = 281
$ws.Range("D66").Value = 279
